$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending=FAPs, Ligand=Wnt5a, Receptor=Fzd3, Target=ECs) ---
$ws.Cells.Item(2,5).Value  = 3
$ws.Cells.Item(2,7).Value  = 4.619088000000001
$ws.Cells.Item(2,8).Value  = 13.857264
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,13).Value = 0.415887
$ws.Cells.Item(2,14).Value = 0.831774
$ws.Cells.Item(2,15).Value = 0.1274765705819998
$ws.Cells.Item(2,16).Value = 0.1197424607254622
$ws.Cells.Item(2,17).Value = 1.921018651056
$ws.Cells.Item(2,18).Value = 11.526111906336
$ws.Cells.Item(2,19).Value = 0.1274765705819998
$ws.Cells.Item(2,20).Value = 0.1197424607254622

# --- Row 3 (Sending=FAPs, Ligand=Wnt5a, Receptor=Fzd3, Target=FAPs) ---
$ws.Cells.Item(3,5).Value  = 3
$ws.Cells.Item(3,7).Value  = 4.619088000000001
$ws.Cells.Item(3,8).Value  = 13.857264
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,13).Value = 0.2747246666666667
$ws.Cells.Item(3,14).Value = 0.824174
$ws.Cells.Item(3,15).Value = 0.08420786983230948
$ws.Cells.Item(3,16).Value = 0.1186483622064973
$ws.Cells.Item(3,17).Value = 1.268977411104
$ws.Cells.Item(3,18).Value = 11.420796699936
$ws.Cells.Item(3,19).Value = 0.08420786983230948
$ws.Cells.Item(3,20).Value = 0.1186483622064973

# --- Row 4 (Sending=FAPs, Ligand=Wnt5a, Receptor=Fzd3, Target=sCs) ---
$ws.Cells.Item(4,5).Value  = 3
$ws.Cells.Item(4,7).Value  = 4.619088000000001
$ws.Cells.Item(4,8).Value  = 13.857264
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.02135966666666667
$ws.Cells.Item(4,14).Value = 0.064079
$ws.Cells.Item(4,15).Value = 0.006547107881326709
$ws.Cells.Item(4,16).Value = 0.009224834078519996
$ws.Cells.Item(4,17).Value = 0.09866217998400001
$ws.Cells.Item(4,18).Value = 0.887959619856
$ws.Cells.Item(4,19).Value = 0.006547107881326709
$ws.Cells.Item(4,20).Value = 0.009224834078519996

# --- Row 5 (NEW): Sending=FAPs, Ligand=Wnt5a, Receptor=Fzd3, Target=Neutro ---
$ws.Cells.Item(5,1).Value  = "FAPs"
$ws.Cells.Item(5,2).Value  = "Wnt5a"
$ws.Cells.Item(5,3).Value  = "Fzd3"
$ws.Cells.Item(5,4).Value  = "Neutro"
$ws.Cells.Item(5,5).Value  = 3
$ws.Cells.Item(5,6).Value  = 1
$ws.Cells.Item(5,7).Value  = 4.619088000000001
$ws.Cells.Item(5,8).Value  = 13.857264
$ws.Cells.Item(5,9).Value  = 1
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.125357
$ws.Cells.Item(5,14).Value = 0.376071
$ws.Cells.Item(5,15).Value = 0.03842409226171471
$ws.Cells.Item(5,16).Value = 0.05413930580600655
$ws.Cells.Item(5,17).Value = 0.5790350144160001
$ws.Cells.Item(5,18).Value = 5.211315129744001
$ws.Cells.Item(5,19).Value = 0.03842409226171471
$ws.Cells.Item(5,20).Value = 0.05413930580600655

# --- Row 6 (NEW): Sending=FAPs, Ligand=Wnt5a, Receptor=Fzd3, Target=M2 ---
$ws.Cells.Item(6,1).Value  = "FAPs"
$ws.Cells.Item(6,2).Value  = "Wnt5a"
$ws.Cells.Item(6,3).Value  = "Fzd3"
$ws.Cells.Item(6,4).Value  = "M2"
$ws.Cells.Item(6,5).Value  = 3
$ws.Cells.Item(6,6).Value  = 1
$ws.Cells.Item(6,7).Value  = 4.619088000000001
$ws.Cells.Item(6,8).Value  = 13.857264
$ws.Cells.Item(6,9).Value  = 1
$ws.Cells.Item(6,10).Value = 1
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.42513
$ws.Cells.Item(6,14).Value = 4.85026
$ws.Cells.Item(6,15).Value = 0.7433443594426494
$ws.Cells.Item(6,16).Value = 0.6982450371835141
$ws.Cells.Item(6,17).Value = 11.20188888144
$ws.Cells.Item(6,18).Value = 67.21133328864001
$ws.Cells.Item(6,19).Value = 0.7433443594426494
$ws.Cells.Item(6,20).Value = 0.6982450371835141
